$d = $word.ActiveDocument

# The document contains a single content control: the auto-generated
# "Table of Contents" field wrapped in a Table of Contents gallery SDT.
# Remove it completely (the TOC heading paragraph plus every TOC entry
# paragraph it contains), leaving the page-break paragraph that used to
# follow it immediately after the date paragraph ("6/30/2019").
$cc = $d.ContentControls.Item(1)
$cc.LockContentControl = $false
$cc.LockContents = $false
$cc.Delete($true)

# The page-break paragraph (now right after "6/30/2019") should carry the
# "_GoBack" bookmark that used to sit further down in the document (right
# after the first "XXX" placeholder in the "Cooperative collision warning
# systems" bullet). Adding a bookmark named "_GoBack" moves/replaces the
# existing one and the remaining heading bookmarks (_Toc...) renumber
# themselves automatically.
$pageBreakPara = $d.Paragraphs.Item(3)
$target = $d.Range($pageBreakPara.Range.Start, $pageBreakPara.Range.Start)
$d.Bookmarks.Add("_GoBack", $target)
